$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"2.5751819521246944E-6"
$ws.Range("A3").Value = [double]"2.3280001641978743E-6"
$ws.Range("F3").Value = 9.0
$ws.Range("A4").Value = [double]"2.4718181634852954E-7"
$ws.Range("F4").Value = 4.0
